# formula: add support for text functions
# Adds a new "Text" worksheet (after "Math and Trig") exercising
# CHAR, CLEAN, CODE, CONCATENATE, EXACT, LEFT, LEN, LOWER, PROPER,
# REPT, T, TRIM, UNICHAR, UNICODE, UPPER.

$wb = $excel.ActiveWorkbook

# --- add the new sheet at the end of the workbook -------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Text"

# --- header row (row 2) -----------------------------------------------
# Ordered (not a hashtable) so shared-string insertion order is
# deterministic: B2..P2 == CHAR, CLEAN, CODE, CONCATENATE, EXACT, LEFT,
# LEN, LOWER, PROPER, REPT, T, TRIM, UNICHAR, UNICODE, UPPER.
$headers = @(
    , @("B2", "CHAR")
    , @("C2", "CLEAN")
    , @("D2", "CODE")
    , @("E2", "CONCATENATE")
    , @("F2", "EXACT")
    , @("G2", "LEFT")
    , @("H2", "LEN")
    , @("I2", "LOWER")
    , @("J2", "PROPER")
    , @("K2", "REPT")
    , @("L2", "T")
    , @("M2", "TRIM")
    , @("N2", "UNICHAR")
    , @("O2", "UNICODE")
    , @("P2", "UPPER")
)

foreach ($pair in $headers) {
    $addr = $pair[0]
    $text = $pair[1]
    $cell = $ws.Range($addr)
    $cell.Value = $text
    $cell.Font.Bold = $true
    if ($addr -ne "C2") {
        $cell.HorizontalAlignment = -4108
    }
}

# --- CHAR ------------------------------------------------------------
$ws.Range("B3").Formula = "=CHAR()"
$ws.Range("B4").Formula = "=CHAR(65)"
$ws.Range("B5").Formula = "=CHAR(33)"
$ws.Range("B6").Formula = "=CHAR(90)"
$ws.Range("B7").Formula = "=CHAR(256)"
$ws.Range("B8").Formula = "=CHAR(-1)"
$ws.Range("B9").Formula = "=CHAR(65.2)"
$ws.Range("B10").Formula = "=CHAR(65.9)"

# --- CLEAN -------------------------------------------------------------
$ws.Range("C3").Formula = "=CLEAN()"
$ws.Range("C4").Formula = "=CLEAN("""")"
$ws.Range("C5").Formula = "=CLEAN(B4)"
$ws.Range("C6").Formula = "=CLEAN(CHAR(9)&""foo""&CHAR(10))"

# --- CODE ----------------------------------------------------------------
$ws.Range("D3").Formula = "=CODE()"
$ws.Range("D4").Formula = "=CODE("""")"
$ws.Range("D5").Formula = "=CODE(B5)"
$ws.Range("D6").Formula = "=CODE(C6)"
$ws.Range("D7").Formula = "=CODE(C5)"

# --- CONCATENATE -----------------------------------------------------------
$ws.Range("E3").Formula = "=CONCATENATE()"
$ws.Range("E4").Formula = "=CONCATENATE(B14:C15)"
$ws.Range("E5").Formula = "=CONCATENATE(B14:C14)"
$ws.Range("E6").Formula = "=CONCATENATE(B14,C14)"

# --- EXACT -------------------------------------------------------------
$ws.Range("F3").Formula = "=EXACT(""a"",""a"")"
$ws.Range("F4").Formula = "=EXACT(""a"",""a "")"
$ws.Range("F5").Formula = "=EXACT(""b"",""b"")"
$ws.Range("F6").Formula = "=EXACT(1,3)"
$ws.Range("F7").Formula = "=EXACT(1,1)"

# --- LEFT ----------------------------------------------------------------
$ws.Range("G3").Formula = "=LEFT(F3)"
$ws.Range("G4").Formula = "=LEFT(D5,2)"
$ws.Range("G5").Formula = "=LEFT(D6,1)"
$ws.Range("G6").Formula = "=LEFT(C6,2)"

# --- LEN -----------------------------------------------------------------
$ws.Range("H3").Formula = "=LEN()"
$ws.Range("H4").Formula = "=LEN(1)"
$ws.Range("H5").Formula = "=LEN(2)"
$ws.Range("H6").Formula = "=LEN(10)"
$ws.Range("H7").Formula = "=LEN(C6)"

# --- LOWER ---------------------------------------------------------------
$ws.Range("I3").Formula = "=LOWER()"
$ws.Range("I4").Formula = "=LOWER(""A"")"
$ws.Range("I5").Formula = "=LOWER(""FOO"")"
$ws.Range("I6").Formula = "=LOWER(""foo"")"

# --- PROPER --------------------------------------------------------------
$ws.Range("J3").Formula = "=PROPER()"
$ws.Range("J4").Formula = "=PROPER("""")"
$ws.Range("J5").Formula = "=PROPER(""foo bar"")"
$ws.Range("J6").Formula = "=PROPER(""foo"")"
$ws.Range("J7").Formula = "=PROPER(""Foo"")"
$ws.Range("J8").Formula = "=PROPER(""foo bar baz   quuz"")"
$ws.Range("J9").Formula = "=PROPER(""foo,bar,baz"")"
$ws.Range("J10").Formula = "=PROPER(""76BudGet"")"

# --- REPT ------------------------------------------------------------------
$ws.Range("K3").Formula = "=REPT()"
$ws.Range("K4").Formula = "=REPT("""")"
$ws.Range("K5").Formula = "=REPT("""",0)"
$ws.Range("K6").Formula = "=REPT(J6,1)"
$ws.Range("K7").Formula = "=REPT(J7,3)"
$ws.Range("K9").Formula = "=REPT(J7,-1)"

# --- T -----------------------------------------------------------------
$ws.Range("L3").Formula = "=T()"
$ws.Range("L4").Formula = "=T(K6)"
$ws.Range("L5").Formula = "=T(H4)"
$ws.Range("L6").Formula = "=T(J3)"

# --- TRIM ----------------------------------------------------------------
$ws.Range("M3").Formula = "=TRIM()"
$ws.Range("M4").Formula = "=TRIM(""A"")"
$ws.Range("M5").Formula = "=TRIM(""A B"")"
$ws.Range("M6").Formula = "=TRIM(""A  B"")"
$ws.Range("M7").Formula = "=TRIM("" A B C   D"")"
$ws.Range("M8").Formula = "=TRIM(""A ""&CHAR(9)&CHAR(9)&"" B  "")"

# --- UNICHAR ---------------------------------------------------------------
$ws.Range("N3").Formula = "=UNICHAR()"
$ws.Range("N4").Formula = "=UNICHAR(32)"
$ws.Range("N5").Formula = "=UNICHAR(66)"

# --- UNICODE ---------------------------------------------------------------
$ws.Range("O3").Formula = "=CODE()"
$ws.Range("O4").Formula = "=UNICODE("""")"
$ws.Range("O5").Formula = "=UNICODE(M5)"
$ws.Range("O6").Formula = "=UNICODE(N6)"
$ws.Range("O7").Formula = "=UNICODE(N5)"

# --- UPPER -----------------------------------------------------------------
$ws.Range("P3").Formula = "=UPPER()"
$ws.Range("P4").Formula = "=UPPER(""A"")"
$ws.Range("P5").Formula = "=UPPER(""a"")"
$ws.Range("P6").Formula = "=UPPER(""foo bar baz"")"

# --- helper data table (B14:C15) ------------------------------------------
$ws.Range("B14").Value = "a"
$ws.Range("C14").Value = "b"
$ws.Range("B15").Value = "c"
$ws.Range("C15").Value = "d"

# --- selection / view state to match the authored file ---------------------
$ws.Range("P7").Select()

Write-Output "Text sheet populated"
